# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "53.782.39"
$ws.Range("E2").Value = "  -11.17%  "
Set-TextValue "D3" "2.329.11"
$ws.Range("E3").Value = "  -19.69%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "441.61"
$ws.Range("E5").Value = "  -16.34%  "
Set-TextValue "D6" "124.93"
$ws.Range("E6").Value = "  -12.90%  "
Set-TextValue "D7" "0.997"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -14.36%  "
Set-TextValue "D9" "2.334.05"
$ws.Range("E9").Value = "  -19.71%  "
Set-TextValue "D10" "5.36"
$ws.Range("E10").Value = "  -11.04%  "
$ws.Range("E12").Value = "  -14.83%  "
$ws.Range("E13").Value = "  -3.33%  "
Set-TextValue "D14" "2.683.33"
$ws.Range("E14").Value = "  -21.19%  "
Set-TextValue "D15" "53.780.28"
$ws.Range("E15").Value = "  -11.18%  "
Set-TextValue "D16" "18.84"
$ws.Range("E17").Value = "  -14.44%  "
Set-TextValue "D18" "2.352.44"
$ws.Range("E18").Value = "  -19.02%  "
Set-TextValue "D19" "3.94"
$ws.Range("E19").Value = "  -21.81%  "
Set-TextValue "D20" "297.49"
$ws.Range("E20").Value = "  -17.73%  "
Set-TextValue "D21" "9.16"
$ws.Range("E21").Value = "  -21.76%  "
Set-TextValue "D22" "0.997"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("E24").Value = "  -19.06%  "
Set-TextValue "D25" "55.66"
$ws.Range("E25").Value = "  -14.02%  "
Set-TextValue "D26" "0.997"
$ws.Range("E26").Value = "  -0.76%  "
Set-TextValue "D27" "0.153"
$ws.Range("E27").Value = "  -15.30%  "
$ws.Range("E28").Value = "  -19.40%  "
Set-TextValue "D29" "6.96"
$ws.Range("E29").Value = "  -11.71%  "
Set-TextValue "D30" "0.996"
$ws.Range("E30").Value = "  -0.29%  "
Set-TextValue "D31" "0.0₃0701"
$ws.Range("E31").Value = "  -17.63%  "
Set-TextValue "D32" "146.37"
$ws.Range("E32").Value = "  -3.83%  "
Set-TextValue "D33" "17.23"
$ws.Range("E33").Value = "  -12.79%  "
$ws.Range("E34").Value = "  -20.12%  "
Set-TextValue "D35" "4.67"
$ws.Range("E35").Value = "  -16.33%  "
Set-TextValue "D36" "3.54"
$ws.Range("E36").Value = "  -19.24%  "
Set-TextValue "D37" "0.828"
$ws.Range("E37").Value = "  -18.00%  "
$ws.Range("E38").Value = "  -16.48%  "
$ws.Range("E39").Value = "  -11.52%  "
Set-TextValue "D40" "0.996"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("E41").Value = "  -0.52%  "
Set-TextValue "D42" "1.942.90"
$ws.Range("E42").Value = "  -15.30%  "
$ws.Range("E43").Value = "  -15.87%  "
Set-TextValue "D44" "1.21"
$ws.Range("E44").Value = "  -18.70%  "
Set-TextValue "D45" "0.0496"
$ws.Range("E45").Value = "  -14.75%  "
Set-TextValue "D46" "0.523"
$ws.Range("E46").Value = "  -19.23%  "
Set-TextValue "D47" "0.0208"
$ws.Range("E47").Value = "  -12.49%  "
$ws.Range("E48").Value = "  -10.18%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "15.96"
$ws.Range("E49").Value = "  -22.07%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "4.02"
$ws.Range("E50").Value = "  -19.55%  "
$ws.Range("E51").Value = "  -3.73%  "
